# "update with new logo and colors"
#
# Metadata sheet changes:
#  - Version 0.1.6 -> 0.1.7
#  - Status active -> draft
#  - Date updated
#  - Contact text updated, and a second Contact row (Bob Milius) is added
#  - A new "Jurisdiction" row (empty value) is inserted after the new Contact row
#  - Everything below shifts down by two rows (Description, Purpose, Copyright,
#    Immutable all move down two rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Simple in-place value updates (rows 1-9 unaffected in position) ---
$ws.Cells.Item(3, 2).Value = "0.1.7"
$ws.Cells.Item(6, 2).Value = "draft"
$ws.Cells.Item(8, 2).Value = "2024-08-23T10:17:11-05:00"

# --- Snapshot the current (pre-shift) values of rows 11-15 before overwriting ---
$a11 = $ws.Cells.Item(11, 1).Value2
$b11 = $ws.Cells.Item(11, 2).Value2
$a12 = $ws.Cells.Item(12, 1).Value2
$b12 = $ws.Cells.Item(12, 2).Value2
$a13 = $ws.Cells.Item(13, 1).Value2
$b13 = $ws.Cells.Item(13, 2).Value2
$a14 = $ws.Cells.Item(14, 1).Value2
$b14 = $ws.Cells.Item(14, 2).Value2
$a15 = $ws.Cells.Item(15, 1).Value2
$b15 = $ws.Cells.Item(15, 2).Value2

# --- Extend the formatted range down to row 16 by copying the formats of row 15 ---
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Shift old rows 11-15 down to 12-16 (bottom-up so nothing is clobbered) ---
$ws.Cells.Item(16, 1).Value = $a15
$ws.Cells.Item(16, 2).Value = $b15
$ws.Cells.Item(15, 1).Value = $a14
$ws.Cells.Item(15, 2).Value = $b14
$ws.Cells.Item(14, 1).Value = $a13
$ws.Cells.Item(14, 2).Value = $b13
$ws.Cells.Item(13, 1).Value = $a12
$ws.Cells.Item(13, 2).Value = $b12
$ws.Cells.Item(12, 1).Value = $a11
$ws.Cells.Item(12, 2).Value = $b11

# --- Row 10: updated Contact text ---
$ws.Cells.Item(10, 2).Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Row 11: new second Contact line ---
$ws.Cells.Item(11, 1).Value = "Contact"
$ws.Cells.Item(11, 2).Value = "Bob Milius (bmilius@nmdp.org)"

# --- Row 12: new Jurisdiction row (value left blank) ---
$ws.Cells.Item(12, 1).Value = "Jurisdiction"
$ws.Cells.Item(12, 2).Value = ""
